$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new data row (2022-01-25) is inserted at row 108,
# pushing all subsequent rows (old 108..209) down by one (new 109..210).
$ws.Rows(108).Insert()

# Populate the newly inserted row 108 with this week's reading.
$ws.Cells.Item(108, 1).Value = 4
$ws.Cells.Item(108, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(108, 3).Value = "Los Lagos"
$ws.Cells.Item(108, 4).Value = 44586
$ws.Cells.Item(108, 5).Value = 10
$ws.Cells.Item(108, 6).Value = 100112003
$ws.Cells.Item(108, 7).Value = "Ajo"
$ws.Cells.Item(108, 8).Value = "Chino"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 240
$ws.Cells.Item(108, 11).Value = 20000
$ws.Cells.Item(108, 12).Value = 21000
$ws.Cells.Item(108, 13).Value = 20500
$ws.Cells.Item(108, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(108, 15).Value = "China"
$ws.Cells.Item(108, 16).Value = 2050
$ws.Cells.Item(108, 17).Value = 10
$ws.Cells.Item(108, 18).Value = "Hortaliza"
